$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values
$ws.Range("B2").Value = 0.0328798185941043
$ws.Range("C2").Value = 0.02902745553985043
$ws.Range("D2").Value = 0.0170043612987258
$ws.Range("E2").Value = 0.0118061167879423
$ws.Range("F2").Value = 0.00716288759837407
$ws.Range("G2").Value = 0.01243154725170473
$ws.Range("H2").Value = 0.01121032945169066
$ws.Range("I2").Value = 0.00600037014603944
$ws.Range("J2").Value = 0.001856398132993877
$ws.Range("K2").Value = 0.001048218029350103
$ws.Range("L2").Value = 0.01091269841269842

# Row 3 values
$ws.Range("B3").Value = 0.0270383024029597
$ws.Range("C3").Value = 0.03108365166631143
$ws.Range("D3").Value = 0.01962557995935889
$ws.Range("E3").Value = 0.0297138470989089
$ws.Range("F3").Value = 0.0252480956350666
$ws.Range("G3").Value = 0.0268442666308431
$ws.Range("H3").Value = 0.0239548245647589
$ws.Range("I3").Value = 0.02295270978019673
$ws.Range("J3").Value = 0.007080252077597466
$ws.Range("K3").Value = 0.009516821407915869
$ws.Range("L3").Value = 0.0334551678646898

# Row 4 values
$ws.Range("B4").Value = 0.020928207938768
$ws.Range("C4").Value = 0.01989440160208194
$ws.Range("D4").Value = 0.01529731160413443
$ws.Range("E4").Value = 0.02246980856607217
$ws.Range("F4").Value = 0.0169934483227596
$ws.Range("G4").Value = 0.0310527569049573
$ws.Range("H4").Value = 0.01356576396703233
$ws.Range("I4").Value = 0.01314400506416592
$ws.Range("J4").Value = 0.003023798163070063
$ws.Range("K4").Value = 0.003406204367085903
$ws.Range("L4").Value = 0.01671263258293975

# Row 5 values (note: G5 does not exist in the sheet, leave it absent)
$ws.Range("B5").Value = 0.0217818530216891
$ws.Range("C5").Value = 0.02318263708627457
$ws.Range("D5").Value = 0.01356396593239173
$ws.Range("E5").Value = 0.02488615365536924
$ws.Range("F5").Value = 0.02986470123147457
$ws.Range("H5").Value = 0.03104114609462904
$ws.Range("I5").Value = 0.01782147620156227
$ws.Range("J5").Value = 0.007902167556687366
$ws.Range("K5").Value = 0.005302742354778724
$ws.Range("L5").Value = 0.01888220638875225
